$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the data (header row + data rows).
$used = $ws.UsedRange
$lastRow = $used.Rows.Count()
$lastCol = $used.Columns.Count()

# Rename header row cells: "..._old" -> "..._FV2410", "..._new" -> "..._FV2504"
for ($col = 1; $col -le $lastCol; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cur = $cell.Value()
    if ($cur -ne $null) {
        if ($cur -match "_old$") {
            $cell.Value = ($cur -replace "_old$", "_FV2410")
        } elseif ($cur -match "_new$") {
            $cell.Value = ($cur -replace "_new$", "_FV2504")
        }
    }
}

# Freeze the header row (split after row 1, freeze panes on).
[void]$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# Turn the data range into an Excel Table ("Table1") with a header row.
$tblRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastRow, $lastCol))
$tbl = $ws.ListObjects.Add(1, $tblRange, $null, 1)
$tbl.Name = "Table1"
